$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert M2, N2, O2 from inline strings into real numeric values
$ws.Range("M2").Value = 61
$ws.Range("N2").Value = 1
$ws.Range("O2").Value = 31

# Add new row 3 of data (kept as text, matching the source diff which stores
# these as inlineStr cells rather than numbers)
$ws.Range("A3").Value = "2025-02-06T17:28"
$ws.Range("B3").Value = "SHIFT1"
$ws.Range("C3").Value = "OK"
$ws.Range("D3").Value = "OK"
$ws.Range("E3").Value = "OK"
$ws.Range("F3").Value = "OK"
$ws.Range("G3").Value = "OK"
$ws.Range("H3").Value = "OK"
$ws.Range("I3").Value = "OK"
$ws.Range("J3").Value = "OK"
$ws.Range("K3").Value = "OK"
$ws.Range("L3").Value = "OK"

$ws.Range("M3:O3").NumberFormat = "@"
$ws.Range("M3").Value = "489"
$ws.Range("N3").Value = "45"
$ws.Range("O3").Value = "54"
